# Add "NumberOfOutOfSchool6to24" (column T) and "NoOfSchool" (column U) to Sheet1.
# Mirrors the author's edit: new header cells in row 1, out-of-school counts for
# most regions in column T (row 12 / Region VIII intentionally left blank), and
# column widths/view state updated to match the new columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers (row 1) ---------------------------------------------------
$ws.Range("T1").Value = "NumberOfOutOfSchool6to24"
$ws.Range("U1").Value = "NoOfSchool"
$ws.Range("T1").Style = $ws.Range("S1").Style
$ws.Range("U1").Style = $ws.Range("S1").Style

# --- Column widths for the two new columns ------------------------------
$ws.Range("T1").ColumnWidth = 25.5703125
$ws.Range("U1").ColumnWidth = 18.140625

# --- Data for column T: NumberOfOutOfSchool6to24 (row 12 / Region VIII is left blank) ---
$outOfSchool = @{
    2  = 4556
    3  = 1806
    4  = 681
    5  = 1261
    6  = 4060
    7  = 5209
    8  = 1189
    9  = 2344
    10 = 2832
    11 = 2776
    13 = 1475
    14 = 1846
    15 = 1862
    16 = 1791
    17 = 1039
    18 = 1511
}

foreach ($row in 2..18) {
    $cell = $ws.Range("T$row")
    if ($outOfSchool.ContainsKey($row)) {
        $cell.Value = $outOfSchool[$row]
    }
    $cell.NumberFormat = "#,##0"
}

# --- View: zoom in a bit and shift the frozen-pane scroll to the new columns ---
$ws.Application.ActiveWindow.Zoom = 115
$ws.Range("U1").Select()
